$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be stored as text so that values like
# "112.00", "0.9999", "25.796.19" keep their exact original formatting
# instead of being auto-converted to numbers (which would strip
# trailing/structural zeros and change the cell type).
$ws.Range("D2:D51").NumberFormat = "@"

# --- Row 25 / Row 26: LidoDAOToken and Monero swapped positions ---
$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"

# --- Price (D) column updates ---
$ws.Range("D2").Value = "25.796.19"
$ws.Range("D3").Value = "1.813.43"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D5").Value = "276.76"
$ws.Range("D6").Value = "0.9995"
$ws.Range("D7").Value = "0.5116"
$ws.Range("D9").Value = "44.71"
$ws.Range("D10").Value = "0.06673"
$ws.Range("D12").Value = "0.8348"
$ws.Range("D13").Value = "0.07827"
$ws.Range("D14").Value = "1.796.58"
$ws.Range("D16").Value = "87.95"
$ws.Range("D17").Value = "0.9991"
$ws.Range("D19").Value = "0.000008023"
$ws.Range("D20").Value = "0.9992"
$ws.Range("D21").Value = "25.870.48"
$ws.Range("D22").Value = "4.732"
$ws.Range("D23").Value = "10.02"
$ws.Range("D24").Value = "6.064"
$ws.Range("D25").Value = "2.206"
$ws.Range("D26").Value = "141.59"
$ws.Range("D27").Value = "1.655"
$ws.Range("D28").Value = "17.04"
$ws.Range("D30").Value = "4.345"
$ws.Range("D31").Value = "4.229"
$ws.Range("D32").Value = "0.08786"
$ws.Range("D33").Value = "0.04883"
$ws.Range("D34").Value = "0.7338"
$ws.Range("D35").Value = "1.141"
$ws.Range("D36").Value = "2.893"
$ws.Range("D38").Value = "3.060"
$ws.Range("D39").Value = "0.5250"
$ws.Range("D40").Value = "0.01858"
$ws.Range("D42").Value = "0.9521"
$ws.Range("D43").Value = "112.00"
$ws.Range("D44").Value = "6.184"
$ws.Range("D45").Value = "8.146"
$ws.Range("D46").Value = "0.9987"
$ws.Range("D47").Value = "0.4581"
$ws.Range("D48").Value = "0.1383"
$ws.Range("D49").Value = "9.293"
$ws.Range("D50").Value = "36.21"
$ws.Range("D51").Value = "1.503"

# --- Volume(1h) (E) column updates ---
$ws.Range("E2").Value = "  -5.11%  "
$ws.Range("E3").Value = "  -4.27%  "
$ws.Range("E5").Value = "  -9.54%  "
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("E7").Value = "  -4.73%  "
$ws.Range("E8").Value = "  -7.10%  "
$ws.Range("E9").Value = "  -1.89%  "
$ws.Range("E10").Value = "  -8.22%  "
$ws.Range("E11").Value = "  -8.52%  "
$ws.Range("E12").Value = "  -6.84%  "
$ws.Range("E13").Value = "  -4.19%  "
$ws.Range("E14").Value = "  -2.39%  "
$ws.Range("E15").Value = "  -4.88%  "
$ws.Range("E16").Value = "  -7.10%  "
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("E18").Value = "  -6.28%  "
$ws.Range("E19").Value = "  -7.10%  "
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("E21").Value = "  -4.26%  "
$ws.Range("E22").Value = "  -5.88%  "
$ws.Range("E23").Value = "  -7.15%  "
$ws.Range("E24").Value = "  -6.14%  "
$ws.Range("E25").Value = "  -3.53%  "
$ws.Range("E26").Value = "  -4.74%  "
$ws.Range("E27").Value = "  -5.60%  "
$ws.Range("E28").Value = "  -6.97%  "
$ws.Range("E29").Value = "  -6.25%  "
$ws.Range("E30").Value = "  -9.72%  "
$ws.Range("E32").Value = "  -4.05%  "
$ws.Range("E33").Value = "  -3.07%  "
$ws.Range("E34").Value = "  -10.31%  "
$ws.Range("E35").Value = "  -6.04%  "
$ws.Range("E36").Value = "  -4.18%  "
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("E38").Value = "  -7.02%  "
$ws.Range("E39").Value = "  -11.69%  "
$ws.Range("E40").Value = "  -6.18%  "
$ws.Range("E41").Value = "  -14.01%  "
$ws.Range("E42").Value = "  -11.37%  "
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("E44").Value = "  -6.57%  "
$ws.Range("E45").Value = "  -11.73%  "
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("E47").Value = "  -9.68%  "
$ws.Range("E48").Value = "  -9.19%  "
$ws.Range("E49").Value = "  -8.52%  "
$ws.Range("E50").Value = "  -4.44%  "
$ws.Range("E51").Value = "  -7.54%  "
